# Trade #27 closed at 2026-02-17 23:57:57 - unknown UNKNOWN +0.000%
#
# Updates the "Summary", "Strategy Status", "All Trades" and
# "MarketMaking" sheets of the live trading results workbook to record
# the newly-closed trade #27 (MarketMaking / DOWN, entry 0.35, exit 0.55,
# closed early for +0.14% / +$0.20, bringing strategy capital to 101.13).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet: refresh the headline stats now that trade #27 closed.
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1501.13   # Current Capital
$wsSummary.Range("B4").Value = 1.13      # Total P&L $
$wsSummary.Range("B5").Value = 0.84      # Total P&L %
$wsSummary.Range("B6").Value = 27        # Total Trades
$wsSummary.Range("B7").Value = 16        # Winning Trades
$wsSummary.Range("B9").Value = 59.26     # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet: update the MarketMaking strategy row (row 6).
# ---------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C6").Value = 101.13     # Capital
$wsStatus.Range("D6").Value = 27         # Trades
$wsStatus.Range("E6").Value = 1.13       # P&L $
$wsStatus.Range("F6").Value = 1.13       # P&L %
$wsStatus.Range("G6").Value = 59.26      # Win Rate %

# ---------------------------------------------------------------------
# Append the newly-closed trade as row 28 on both the "All Trades" and
# "MarketMaking" logs (they carry an identical copy of the trade log).
# ---------------------------------------------------------------------
function Add-TradeRow($ws) {
    $row = 28

    # Column B (Date) must stay plain text ("2026-02-17"), not get
    # auto-coerced into an Excel date serial.
    $ws.Cells.Item($row, 2).NumberFormat = "@"

    $ws.Cells.Item($row, 1).Value  = 27
    $ws.Cells.Item($row, 2).Value  = "2026-02-17"
    $ws.Cells.Item($row, 3).Value  = "23:57:51"
    $ws.Cells.Item($row, 4).Value  = "MarketMaking"
    $ws.Cells.Item($row, 5).Value  = "DOWN"
    $ws.Cells.Item($row, 6).Value  = 0.35
    $ws.Cells.Item($row, 7).Value  = 0.55
    $ws.Cells.Item($row, 8).Value  = "CLOSED"
    $ws.Cells.Item($row, 9).Value  = 57.1429
    $ws.Cells.Item($row, 10).Value = 0.2
    $ws.Cells.Item($row, 11).Value = 101.13
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($row, 16).Value = "early_exit"
    $ws.Cells.Item($row, 17).Value = 0.14
}

Add-TradeRow($wb.Worksheets.Item("All Trades"))
Add-TradeRow($wb.Worksheets.Item("MarketMaking"))
